$d = $word.ActiveDocument

# --- Change 2 (done first) ---
# The document currently has exactly one "_GoBack" bookmark, sitting alone
# in its own empty paragraph near the end of the document. Remove it now,
# before Change 1 below introduces a new "_GoBack" bookmark elsewhere -
# doing it in this order means the lookup by name is never ambiguous.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Change 1 ---
# Paragraph "nộp đơn: nhập ngày nộp, tải file copy lên -> ấn ok thì chuyển
# trạng thái và gửi email cho khách hàng (đính kèm 2 file: advice filing +
# bản tải)" needs its first four runs collapsed into:
#   "nộp đơn: " + (bookmark _GoBack, moved here) + "khách hàng (đính kèm 2 file: "
# while leaving the final run "advice filing + bản tải)" as-is.
#
# We locate the paragraph via Find (unique anchor text), then rewrite its
# range using InsertXML so we get full control of the resulting run
# boundaries (plain Find/Replace normalizes all touched runs into one,
# which would also swallow the trailing "advice filing + bản tải)" run).

$anchor = "nộp đơn: nhập ngày nộp"
$found = $d.Content.Find.Execute($anchor, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$para = $d.Paragraphs.Item(16)
$rng = $para.Range

$fragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
  '<w:body>' +
  '<w:p w14:paraId="3D1FAC14" w14:textId="22EFBD8F" w:rsidR="00BE2F9B" w:rsidRDefault="00BE2F9B" w:rsidP="00BE2F9B">' +
  '<w:pPr><w:pStyle w:val="oancuaDanhsach"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">nộp đơn: </w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:t xml:space="preserve">khách hàng (đính kèm 2 file: </w:t></w:r>' +
  '<w:r w:rsidR="00253623"><w:t>advice filing + bản tải)</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($fragment)
